# Adding new feature, deleting the student row for the given roll no
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the student row for roll no 107 (Anikit) - entire row shift up
$ws.Rows.Item(5).Delete()

# Delete the rows for roll no 112 (Prakrity Maddheshiya) and 120 (shub)
# After the previous delete, these are now rows 7 and 8
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# Reset leave counters for the remaining students / update records
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

$ws.Range("C4").Value = "nikhilside72@gmail.com"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
